# Apply the "Leetcode problem two sum and others" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width (new column used for the Description header) ---
$ws.Columns.Item(2).ColumnWidth = 110.1640625

# Keep a handle on the existing hyperlink-styled cell so we can clone its
# style onto the freshly added hyperlink cells further down.
$hyperlinkStyle = $ws.Range("A5").Style
$titleStyle = $ws.Range("A1").Style

# --- Row 1: add "Description" header in B1, matching the A1 title style ---
$ws.Range("B1").Value = "Description"
$ws.Range("B1").Style = $titleStyle

# --- "Strings" category (A9) renamed to "String Manipulation" ---
$ws.Range("A9").Value = "String Manipulation"

# --- New rows under "String Manipulation" (A11, A12) ---
$ws.Range("A11").Value = "https://leetcode.com/problems/number-of-senior-citizens/"
$ws.Range("A11").Style = $hyperlinkStyle
$ws.Hyperlinks.Add($ws.Range("A11"), "https://leetcode.com/problems/number-of-senior-citizens/")
$ws.Range("A11").Style = $hyperlinkStyle

$ws.Range("A12").Value = "https://leetcode.com/problems/score-of-a-string/"
$ws.Range("A12").Style = $hyperlinkStyle
$ws.Hyperlinks.Add($ws.Range("A12"), "https://leetcode.com/problems/score-of-a-string/")
$ws.Range("A12").Style = $hyperlinkStyle

# --- New rows under "Hashmap & Hashset" (A15, A16) ---
$ws.Range("A15").Value = "https://leetcode.com/problems/two-sum/description/"
$ws.Range("A15").Style = $hyperlinkStyle
$ws.Hyperlinks.Add($ws.Range("A15"), "https://leetcode.com/problems/two-sum/description/")
$ws.Range("A15").Style = $hyperlinkStyle

$ws.Range("A16").Value = "https://leetcode.com/problems/valid-anagram/"
$ws.Range("A16").Style = $hyperlinkStyle
$ws.Hyperlinks.Add($ws.Range("A16"), "https://leetcode.com/problems/valid-anagram/")
$ws.Range("A16").Style = $hyperlinkStyle

# --- New "Sliding Window" section (A22, A23, B23) ---
$ws.Range("A22").Value = "Sliding Window"
$ws.Range("A22").Style = $titleStyle

$ws.Range("A23").Value = "https://leetcode.com/problems/max-consecutive-ones/"
$ws.Range("A23").Style = $hyperlinkStyle
$ws.Hyperlinks.Add($ws.Range("A23"), "https://leetcode.com/problems/max-consecutive-ones/")
$ws.Range("A23").Style = $hyperlinkStyle

$ws.Range("B23").Value = "keep count of 1s and increase the window when we find subs ones …..reset window when no subs ones"

# --- Selection / view state ---
$ws.Range("A26").Select()
